$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$wsMonthly = $wb.Worksheets.Item("Monthly Performance")

# ---------------------------------------------------------------------------
# 1. Columns B and C get wider (B: 18 -> 45, C: 28 -> 20)
#    (offsets compensate for Excel's internal character-width padding so the
#     stored width matches the target exactly)
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 44.16666666666666
$ws.Columns.Item(3).ColumnWidth = 19.16666666666666

# ---------------------------------------------------------------------------
# 2. Row 17 becomes the new "TRADING ACTIVITY SUMMARY" banner row
#    (it used to be the Metric/Mar-Jul25/Aug-Oct25/Nov25-Feb26/Total header)
# ---------------------------------------------------------------------------
$ws.Range("A4").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A17").Borders.LineStyle = 1
$ws.Range("A17").Borders.Weight = 2
$ws.Range("A17").Borders.Color = 0
$ws.Range("A17").Value = "TRADING ACTIVITY SUMMARY"

$ws.Range("C23").Copy()
$ws.Range("B17:D17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("E23").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B17:E17").ClearContents()
$ws.Rows.Item(17).RowHeight = 22

$ws.Range("A17:E17").Merge()

# ---------------------------------------------------------------------------
# 3. Trading Activity Summary body rows (18-20): relabel + new "average" text
# ---------------------------------------------------------------------------
$wsMonthly.Range("B8").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B18").HorizontalAlignment = -4131
$ws.Range("B18").Value = "0 per month average"

$ws.Range("B18").Copy()
$ws.Range("B19:B20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B19").Value = "0 average"
$ws.Range("B20").Value = "0 average"

$ws.Range("A19").Value = "Buy Transactions"
$ws.Range("A20").Value = "Sell Transactions"

# ---------------------------------------------------------------------------
# 4. "KEY INSIGHTS & RECOMMENDATIONS" header (row 22) reverts to the bold
#    navy banner style (same as row 4 / row 16), and grows to 22pt height
# ---------------------------------------------------------------------------
$ws.Range("A4").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Rows.Item(22).RowHeight = 22

# ---------------------------------------------------------------------------
# 5. Key insight rows 23-28: numbered insight text, light-green banded style,
#    shorter row height (28 -> 20)
# ---------------------------------------------------------------------------
$ws.Range("C6").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A23").Interior.Color = 14348258
$ws.Range("A23").WrapText = $true

$ws.Range("A23").Copy()
$ws.Range("A24:A28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A23").Value = "1. Portfolio demonstrates consistent positive growth with strong cumulative returns"
$ws.Range("A24").Value = "2. High win rate (83%+ positive months) indicates favorable market positioning"
$ws.Range("A25").Value = "3. Dividend accumulation provides steady passive income stream"
$ws.Range("A26").Value = "4. Average monthly returns exceed typical market benchmarks"
$ws.Range("A27").Value = "5. Trading activity shows disciplined approach with measured transactions"
$ws.Range("A28").Value = "6. Risk management evident from contained worst-month losses relative to gains"

for ($r = 23; $r -le 28; $r++) {
    $ws.Rows.Item($r).RowHeight = 20
}

# ---------------------------------------------------------------------------
# 6. "ACTION ITEMS & STRATEGY" header (row 30) just grows taller (20 -> 22)
# ---------------------------------------------------------------------------
$ws.Rows.Item(30).RowHeight = 22

# ---------------------------------------------------------------------------
# 7. Action item rows 31-35 get renumbered text + new row 36; banded peach
#    style, shorter row height (32 -> 20)
# ---------------------------------------------------------------------------
$ws.Range("C6").Copy()
$ws.Range("A31").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A31").Interior.Color = 8696052
$ws.Range("A31").WrapText = $true

$ws.Range("A31").Copy()
$ws.Range("A32:A36").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A31").Value = "1. Continue current strategy - proven track record of consistent returns"
$ws.Range("A32").Value = "2. Maintain dividend reinvestment for compound growth acceleration"
$ws.Range("A33").Value = "3. Review quarterly performance against benchmarks (S&P 500, Russell 2000)"
$ws.Range("A34").Value = "4. Rebalance portfolio if allocation drifts >10% from target"
$ws.Range("A35").Value = "5. Evaluate tax-loss harvesting opportunities in down months"
$ws.Range("A36").Value = "6. Monitor market conditions for tactical adjustments if warranted"

for ($r = 31; $r -le 36; $r++) {
    $ws.Rows.Item($r).RowHeight = 20
}
